$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 799
$ws.Range("I43").Value = 600
$ws.Range("K43").Value = 600
$ws.Range("M43").Value = -531

$ws.Range("H98").Value = 1111
$ws.Range("I98").Value = 639.1667
$ws.Range("K98").Value = 639.1667
$ws.Range("M98").Value = 858.8333

$ws.Range("H122").Value = 1111
$ws.Range("I122").Value = 639.1667
$ws.Range("K122").Value = 1917.5001
$ws.Range("M122").Value = 532.4999

$ws.Range("H129").Value = 1396
$ws.Range("I129").Value = 810
$ws.Range("J129").Value = 2735.4285
$ws.Range("K129").Value = 2430
$ws.Range("L129").Value = 8206.2855
$ws.Range("M129").Value = 2570
$ws.Range("N129").Value = -18206.2855

$ws.Range("H132").Value = 2150
$ws.Range("I132").Value = 2204.8096
$ws.Range("K132").Value = 6614.4288
$ws.Range("M132").Value = -4084.4288

$ws.Range("H137").Value = 3240.2354
$ws.Range("I137").Value = 1682.3334
$ws.Range("K137").Value = 5047.0002
$ws.Range("M137").Value = -2497.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2827.8333
$ws.Range("I61").Value = 2827.8333
$ws.Range("K61").Value = 2827.8333
$ws.Range("M61").Value = -2615.8333

$ws.Range("H74").Value = 6358
$ws.Range("I74").Value = 4898.4
$ws.Range("J74").Value = 10007
$ws.Range("K74").Value = 4898.4
$ws.Range("L74").Value = 10007
$ws.Range("M74").Value = -4024.4
$ws.Range("N74").Value = -11755

$ws.Range("H77").Value = 6358
$ws.Range("I77").Value = 4898.4
$ws.Range("J77").Value = 10007
$ws.Range("K77").Value = 24492
$ws.Range("L77").Value = 50035
$ws.Range("M77").Value = -20124
$ws.Range("N77").Value = -58771

$ws.Range("H110").Value = 2173.6365
$ws.Range("I110").Value = 2102.2222
$ws.Range("K110").Value = 2102.2222
$ws.Range("M110").Value = -57.22220000000016

$ws.Range("H122").Value = 1300
$ws.Range("I122").Value = 1300
$ws.Range("K122").Value = 3900
$ws.Range("M122").Value = -1450

$ws.Range("H132").Value = 2084.6365
$ws.Range("I132").Value = 1491.5
$ws.Range("J132").Value = 3666.3333
$ws.Range("K132").Value = 4474.5
$ws.Range("L132").Value = 10998.9999
$ws.Range("M132").Value = -1944.5
$ws.Range("N132").Value = -16058.9999

$ws.Range("H136").Value = 2827.8333
$ws.Range("I136").Value = 2827.8333
$ws.Range("K136").Value = 8483.499899999999
$ws.Range("M136").Value = -5933.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 337.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 337.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 337.5
$ws.Range("N7").Value = -563.5
$ws.Range("M7").ClearContents()

$ws.Range("H20").Value = 6998.5
$ws.Range("I20").Value = 6998.6665
$ws.Range("K20").Value = 6998.6665
$ws.Range("M20").Value = -6751.6665

$ws.Range("H94").Value = 3490.3635
$ws.Range("I94").Value = 3339.4
$ws.Range("K94").Value = 3339.4
$ws.Range("M94").Value = -2888.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I22").Value = 75.5
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 75.5
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 274.5
$ws.Range("N22").Value = -900

$ws.Range("H31").Value = 3449.05
$ws.Range("J31").Value = 3955.3
$ws.Range("L31").Value = 3955.3
$ws.Range("N31").Value = -4545.3

$ws.Range("H34").Value = 3449.05
$ws.Range("J34").Value = 3955.3
$ws.Range("L34").Value = 3955.3
$ws.Range("N34").Value = -4359.3

$ws.Range("H58").Value = 2508
$ws.Range("J58").Value = 2244
$ws.Range("L58").Value = 2244
$ws.Range("N58").Value = -2650

$ws.Range("H86").Value = 7000
$ws.Range("I86").Value = 7000
$ws.Range("K86").Value = 7000
$ws.Range("M86").Value = -5877

$ws.Range("H89").Value = 7000
$ws.Range("I89").Value = 7000
$ws.Range("K89").Value = 35000
$ws.Range("M89").Value = -29384

$ws.Range("H97").Value = 11500
$ws.Range("I97").Value = 5000
$ws.Range("J97").Value = 18000
$ws.Range("K97").Value = 5000
$ws.Range("L97").Value = 18000
$ws.Range("M97").Value = -4009
$ws.Range("N97").Value = -19982

$ws.Range("H109").Value = 69990
$ws.Range("J109").Value = 69990
$ws.Range("L109").Value = 69990
$ws.Range("N109").Value = -72070

$ws.Range("H132").Value = 981.6875
$ws.Range("I132").Value = 981.6875
$ws.Range("K132").Value = 2945.0625
$ws.Range("M132").Value = -415.0625

$ws.Range("H136").Value = 2508
$ws.Range("J136").Value = 2244
$ws.Range("L136").Value = 6732
$ws.Range("N136").Value = -11832

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 9812
$ws.Range("I62").Value = 9812
$ws.Range("K62").Value = 29436
$ws.Range("M62").Value = -28750

$ws.Range("H65").Value = 9812
$ws.Range("I65").Value = 9812
$ws.Range("K65").Value = 88308
$ws.Range("M65").Value = -84876

$ws.Range("H106").Value = 5024.1665
$ws.Range("J106").Value = 5024.1665
$ws.Range("L106").Value = 15072.4995
$ws.Range("N106").Value = -16964.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 17861536
$ws.Range("I122").Value = 25003048
$ws.Range("J122").Value = 7754
$ws.Range("K122").Value = 75009144
$ws.Range("L122").Value = 23262
$ws.Range("M122").Value = -75006694
$ws.Range("N122").Value = -28162

$ws.Range("H126").Value = 7033.3335
$ws.Range("I126").Value = 8250
$ws.Range("J126").Value = 4600
$ws.Range("K126").Value = 24750
$ws.Range("L126").Value = 13800
$ws.Range("M126").Value = -22280
$ws.Range("N126").Value = -18740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2007.0834
$ws.Range("I22").Value = 2010.2222
$ws.Range("J22").Value = 1997.6666
$ws.Range("K22").Value = 2010.2222
$ws.Range("L22").Value = 1997.6666
$ws.Range("M22").Value = -1715.2222
$ws.Range("N22").Value = -2587.6666

$ws.Range("H27").Value = 2007.0834
$ws.Range("I27").Value = 2010.2222
$ws.Range("J27").Value = 1997.6666
$ws.Range("K27").Value = 2010.2222
$ws.Range("L27").Value = 1997.6666
$ws.Range("M27").Value = -1903.2222
$ws.Range("N27").Value = -2211.6666

$ws.Range("H40").Value = 2498.5
$ws.Range("I40").Value = 2498.5
$ws.Range("K40").Value = 2498.5
$ws.Range("M40").Value = -2362.5

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

$ws.Range("H61").Value = 675.75
$ws.Range("I61").Value = 675.75
$ws.Range("K61").Value = 675.75
$ws.Range("M61").Value = -473.75

$ws.Range("H62").Value = 30300
$ws.Range("J62").Value = 30300
$ws.Range("L62").Value = 30300
$ws.Range("N62").Value = -31548

$ws.Range("H65").Value = 30300
$ws.Range("J65").Value = 30300
$ws.Range("L65").Value = 90900
$ws.Range("N65").Value = -97140

$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 1000
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 248
$ws.Range("N93").Value = -3496

$ws.Range("H113").Value = 675.75
$ws.Range("I113").Value = 675.75
$ws.Range("K113").Value = 675.75
$ws.Range("M113").Value = 1494.25

$ws.Range("H123").Value = 79997
$ws.Range("J123").Value = 79997
$ws.Range("L123").Value = 79997
$ws.Range("N123").Value = -89797

$ws.Range("H132").Value = 6599.857
$ws.Range("I132").Value = 2349.5
$ws.Range("K132").Value = 7048.5
$ws.Range("M132").Value = -4518.5

$ws.Range("H136").Value = 3050.4443
$ws.Range("I136").Value = 2931.875
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 8795.625
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -6245.625
$ws.Range("N136").Value = -17097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 15000
$ws.Range("J33").Value = 10000
$ws.Range("L33").Value = 10000
$ws.Range("N33").Value = -10500

$ws.Range("H36").Value = 15000
$ws.Range("J36").Value = 10000
$ws.Range("L36").Value = 10000
$ws.Range("N36").Value = -10500

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H68").Value = 44450
$ws.Range("J68").Value = 44450
$ws.Range("L68").Value = 44450
$ws.Range("N68").Value = -46072

$ws.Range("H71").Value = 44450
$ws.Range("J71").Value = 44450
$ws.Range("L71").Value = 133350
$ws.Range("N71").Value = -141462

$ws.Range("H113").Value = 5749.3687
$ws.Range("I113").Value = 8618.083000000001
$ws.Range("K113").Value = 25854.249
$ws.Range("M113").Value = -23684.249

$ws.Range("H132").Value = 2185.2144
$ws.Range("I132").Value = 2163.182
$ws.Range("J132").Value = 2266
$ws.Range("K132").Value = 6489.545999999999
$ws.Range("L132").Value = 6798
$ws.Range("M132").Value = -3959.545999999999
$ws.Range("N132").Value = -11858
